$wb = $excel.ActiveWorkbook

# --- Battle 2: nudge the saved selection (L3 -> O2), without leaving it the active tab ---
$ws2 = $wb.Worksheets.Item("Battle 2")
$ws2.Range("O2").Select()

# --- Battle 4: built as a copy of Battle 2 (same column widths/layout), placed after Battle 3 ---
$ws3 = $wb.Worksheets.Item("Battle 3")
$ws2.Copy($null, $ws3)

$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "Battle 4"

# Row 2: rename the player from MUTE to Zappo, command Nail -> Fire, target Eagle -> Goblin
$ws4.Range("B2").Value = "Zappo"
$ws4.Range("L2").Value = "Fire"
$ws4.Range("M2").Value = "Goblin"

# Row 3: the enemy becomes a Goblin with 8 lives
$ws4.Range("B3").Value = "Goblin"
$ws4.Range("D3").Value = 8

# Leave the selection on the next empty row, like the other freshly-edited sheets
$ws4.Rows.Item(4).Select()
